$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ordering of id (B) / speaker_variant (C) pairs for rows 2-11,
# with is_prefered (D) cleared for all these rows.
$rows = @(
    @{ Id = "#diadumen";    Variant = "Diadumen" },
    @{ Id = "#darida";      Variant = "Darida" },
    @{ Id = "#adrastus";    Variant = "Adrastus" },
    @{ Id = "#sabellus";    Variant = "Sabellus" },
    @{ Id = "#fuluius";     Variant = "Fuluius" },
    @{ Id = "#placidus";    Variant = "Placidus" },
    @{ Id = "#liuia";       Variant = "Liuia" },
    @{ Id = "#liuia-geest"; Variant = "Liuia geest" },
    @{ Id = "#labinus";     Variant = "Labinus" },
    @{ Id = "#keyser";      Variant = "Keyser" }
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 2).Value = $row.Id
    $ws.Cells.Item($r, 3).Value = $row.Variant
    $ws.Cells.Item($r, 4).Value = ""
    $r++
}
